# Auto-generated: update Price (D) and Volume(1h) (E) columns for cryptos sheet
# per the commit diff (GitHub Actions refresh of cryptos list).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    # Force the literal text into the cell (avoids Excel auto-converting
    # numeric-looking strings like "213.98" or "0.0890" into numbers),
    # then strip the temporary text-format flag so the cell's style index
    # stays identical to the original (no numFmt/quotePrefix left behind).
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue ($ws.Range('D2')) '27.151.63'
Set-TextValue ($ws.Range('E2')) '  +0.42%  '
Set-TextValue ($ws.Range('D3')) '1.677.93'
Set-TextValue ($ws.Range('E3')) '  -0.18%  '
Set-TextValue ($ws.Range('E4')) '  +0.13%  '
Set-TextValue ($ws.Range('D5')) '213.98'
Set-TextValue ($ws.Range('E5')) '  -0.98%  '
Set-TextValue ($ws.Range('E6')) '  -0.22%  '
Set-TextValue ($ws.Range('E7')) '  +0.12%  '
Set-TextValue ($ws.Range('D8')) '22.64'
Set-TextValue ($ws.Range('E8')) '  +5.43%  '
Set-TextValue ($ws.Range('E9')) '  +1.78%  '
Set-TextValue ($ws.Range('D10')) '0.0621'
Set-TextValue ($ws.Range('E10')) '  -0.34%  '
Set-TextValue ($ws.Range('D11')) '0.0890'
Set-TextValue ($ws.Range('E11')) '  +0.04%  '
Set-TextValue ($ws.Range('D12')) '1.914.89'
Set-TextValue ($ws.Range('E12')) '  -0.16%  '
Set-TextValue ($ws.Range('D13')) '1.679.24'
Set-TextValue ($ws.Range('E13')) '  -0.27%  '
Set-TextValue ($ws.Range('E14')) '  +1.84%  '
Set-TextValue ($ws.Range('D15')) '0.551'
Set-TextValue ($ws.Range('E15')) '  +3.39%  '
Set-TextValue ($ws.Range('D16')) '66.50'
Set-TextValue ($ws.Range('E16')) '  -0.09%  '
Set-TextValue ($ws.Range('D17')) '27.120.97'
Set-TextValue ($ws.Range('E17')) '  +0.29%  '
Set-TextValue ($ws.Range('D18')) '235.56'
Set-TextValue ($ws.Range('E18')) '  +0.05%  '
Set-TextValue ($ws.Range('D19')) '7.86'
Set-TextValue ($ws.Range('E19')) '  -3.95%  '
Set-TextValue ($ws.Range('D20')) '0.0₃0739'
Set-TextValue ($ws.Range('E20')) '  +0.04%  '
Set-TextValue ($ws.Range('E21')) '  +0.15%  '
Set-TextValue ($ws.Range('D22')) '4.52'
Set-TextValue ($ws.Range('E22')) '  +1.33%  '
Set-TextValue ($ws.Range('D23')) '9.53'
Set-TextValue ($ws.Range('E23')) '  +2.63%  '
Set-TextValue ($ws.Range('D24')) '2.09'
Set-TextValue ($ws.Range('E24')) '  -1.49%  '
Set-TextValue ($ws.Range('D25')) '146.54'
Set-TextValue ($ws.Range('E25')) '  +0.02%  '
Set-TextValue ($ws.Range('D26')) '7.39'
Set-TextValue ($ws.Range('E26')) '  +1.97%  '
Set-TextValue ($ws.Range('D27')) '16.32'
Set-TextValue ($ws.Range('E27')) '  -0.73%  '
Set-TextValue ($ws.Range('E28')) '  -0.14%  '
Set-TextValue ($ws.Range('E29')) '  +0.23%  '
Set-TextValue ($ws.Range('D30')) '0.0501'
Set-TextValue ($ws.Range('E30')) '  +0.46%  '
Set-TextValue ($ws.Range('E31')) '  -0.49%  '
Set-TextValue ($ws.Range('E32')) '  -0.04%  '
Set-TextValue ($ws.Range('D33')) '1.539.29'
Set-TextValue ($ws.Range('E33')) '  +0.34%  '
Set-TextValue ($ws.Range('D34')) '3.23'
Set-TextValue ($ws.Range('E34')) '  +1.63%  '
Set-TextValue ($ws.Range('D35')) '1.66'
Set-TextValue ($ws.Range('E35')) '  -3.21%  '
Set-TextValue ($ws.Range('D36')) '0.604'
Set-TextValue ($ws.Range('E36')) '  +2.37%  '
Set-TextValue ($ws.Range('D37')) '0.941'
Set-TextValue ($ws.Range('E37')) '  +1.80%  '
Set-TextValue ($ws.Range('D38')) '2.39'
Set-TextValue ($ws.Range('E38')) '  -0.20%  '
Set-TextValue ($ws.Range('E39')) '  -1.66%  '
Set-TextValue ($ws.Range('E40')) '  +3.39%  '
Set-TextValue ($ws.Range('E41')) '  +1.12%  '
Set-TextValue ($ws.Range('D42')) '69.19'
Set-TextValue ($ws.Range('E42')) '  +1.95%  '
Set-TextValue ($ws.Range('E43')) '  +0.11%  '
Set-TextValue ($ws.Range('E44')) '  -0.14%  '
Set-TextValue ($ws.Range('D45')) '1.822.40'
Set-TextValue ($ws.Range('E45')) '  +0.00%  '
Set-TextValue ($ws.Range('D46')) '0.788'
Set-TextValue ($ws.Range('E46')) '  +0.68%  '
Set-TextValue ($ws.Range('D47')) '89.58'
Set-TextValue ($ws.Range('E47')) '  -0.58%  '
Set-TextValue ($ws.Range('E48')) '  +3.68%  '
Set-TextValue ($ws.Range('E49')) '  +6.07%  '
Set-TextValue ($ws.Range('E50')) '  +2.82%  '
Set-TextValue ($ws.Range('D51')) '0.104'
Set-TextValue ($ws.Range('E51')) '  -0.19%  '
